$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# header style already used by B1:H1 (copy format from H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I2:I85 and J2:J85
$iValues = @(3,7,7,6,8,7,8,8,8,8,6,7,6,8,6,7,4,9,9,9,8,8,8,8,9,7,7,8,7,6,8,6,8,6,6,8,7,8,4,7,7,5,6,6,6,6,4,7,4,4,8,6,7,7,7,4,7,9,6,7,9,8,10,7,7,6,7,7,7,6,7,10,9,7,6,6,7,7,6,7,7,6,4,3)
$jValues = @(4,7,8,7,8,8,8,8,8,8,6,7,7,8,6,7,5,9,9,9,8,8,8,8,9,7,8,8,7,6,8,7,8,7,6,8,7,8,5,7,7,6,6,6,7,7,4,7,5,5,8,6,7,7,7,4,7,9,6,7,9,9,10,7,7,6,7,7,7,6,7,10,9,7,6,6,7,7,6,7,7,6,4,3)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}

Write-Output "done"
